$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pythonCode")

$LF = [char]10

# ---------------------------------------------------------------------------
# 1. Remove the trailing rows (old rows 10, 11, 12) - the sheet shrinks from
#    12 rows (A1:B12) down to 9 rows (A1:B9).
# ---------------------------------------------------------------------------
$ws.Rows("10:12").Delete()

# ---------------------------------------------------------------------------
# 2. Rewrite the python-snippet text in column A (rows 2-9) - every literal
#    "\b" token becomes "\xc" (content-only change, same cell style/shape).
# ---------------------------------------------------------------------------
$searchCode = 'def search(input_list, num):' + $LF + `
  'if(num in input_list):' + $LF + `
  'print("Element Found")' + $LF + `
  '\xc' + $LF + `
  '\xc' + $LF + `
  'else:' + $LF + `
  'print("Not Found")' + $LF + `
  '\xc' + $LF + `
  '\xc' + $LF + `
  '\xc' + $LF + `
  '\xc' + $LF + `
  'search([12, 23, 45, 67, 6, 90] , 12)'

$maxOnesCode = 'def findMaxConsecutiveOnes(nums) :' + $LF + `
  'count = 0' + $LF + `
  'result = 0' + $LF + `
  'for i in range(0, len(nums)):' + $LF + `
  'if (nums[i] == 0):' + $LF + `
  'count = 0' + $LF + `
  '\xc' + $LF + `
  '\xc' + $LF + `
  'else:' + $LF + `
  'count+= 1' + $LF + `
  '\xc' + $LF + `
  '\xc' + $LF + `
  'result = max(result, count)' + $LF + `
  '\xc' + $LF + `
  '\xc' + $LF + `
  'print(result)' + $LF + `
  '\xc' + $LF + `
  '\xc' + $LF + `
  'findMaxConsecutiveOnes([1,0,1,1,0,1])'

$findNumbersCode = 'def findNumbers(nums):' + $LF + `
  'c=0' + $LF + `
  'for i in nums:' + $LF + `
  'j=str(i)' + $LF + `
  'x=len(j)' + $LF + `
  'if x%2==0:' + $LF + `
  'c=c+1' + $LF + `
  '\xc' + $LF + `
  '\xc' + $LF + `
  '\xc' + $LF + `
  '\xc' + $LF + `
  'print c' + $LF + `
  'return c' + $LF + `
  'findNumbers([12,345,2,6,7896])'

$sortedSquaresCode = 'def sortedSquares(nums):' + $LF + `
  'squares_list = []' + $LF + `
  'for i in range(0, len(nums)):' + $LF + `
  'square = nums[i] * nums[i];' + $LF + `
  'squares_list.append(square)' + $LF + `
  '\xc' + $LF + `
  '\xc' + $LF + `
  'sorted_squares_list = sorted(squares_list)' + $LF + `
  'print sorted_squares_list;' + $LF + `
  'return sorted_squares_list;' + $LF + `
  'sortedSquares([-7,-3,2,3,11])'

$ws.Cells.Item(2,1).Value2 = $searchCode
$ws.Cells.Item(3,1).Value2 = $searchCode
$ws.Cells.Item(4,1).Value2 = $maxOnesCode
$ws.Cells.Item(5,1).Value2 = $maxOnesCode
$ws.Cells.Item(6,1).Value2 = $findNumbersCode
$ws.Cells.Item(7,1).Value2 = $findNumbersCode
$ws.Cells.Item(8,1).Value2 = $sortedSquaresCode
$ws.Cells.Item(9,1).Value2 = $sortedSquaresCode

# ---------------------------------------------------------------------------
# 3. Row heights: rows 2/3/9 = 180, rows 4/5 = 285, rows 6/7 = 210, row 8 = 180
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 180
$ws.Rows.Item(3).RowHeight = 180
$ws.Rows.Item(4).RowHeight = 285
$ws.Rows.Item(5).RowHeight = 285
$ws.Rows.Item(6).RowHeight = 210
$ws.Rows.Item(7).RowHeight = 210
$ws.Rows.Item(8).RowHeight = 180
$ws.Rows.Item(9).RowHeight = 180

# ---------------------------------------------------------------------------
# 4. Rewrite column B (the "Result" column) with the new set of outcomes,
#    applying formatting BEFORE the value so look-alike numbers (e.g. "2")
#    are stored as text, not auto-converted to a numeric literal.
#
#    - "Element Found" (B2): plain Aptos Narrow text, no wrap.
#    - "2" / "[4, 9, 9, 49, 121]" (B4, B6, B8): Aptos Narrow, text ("@") format.
#    - " Some Tests failed..." / " No tests were collected" (B3, B5, B9):
#      Arial Unicode MS 12pt black, vertically centered.
#    - "Error occurred during submission" (B7): Consolas 10pt light-gray.
# ---------------------------------------------------------------------------
$someTestsFailed = " Some Tests failed. Please review code"
$errorOccurred = "Error occurred during submission"
$noTestsCollected = " No tests were collected"

$b2 = $ws.Cells.Item(2,2)
$b2.Font.Name = "Aptos Narrow"
$b2.Font.Size = 11
$b2.WrapText = $false
$b2.Value2 = "Element Found"

foreach ($pair in @(@(4,"2"), @(6,"2"), @(8,"[4, 9, 9, 49, 121]"))) {
    $r = $pair[0]
    $val = $pair[1]
    $cell = $ws.Cells.Item($r,2)
    $cell.Font.Name = "Aptos Narrow"
    $cell.Font.Size = 11
    $cell.NumberFormat = "@"
    $cell.WrapText = $false
    $cell.Value2 = $val
}

foreach ($pair in @(@(3,$someTestsFailed), @(5,$someTestsFailed), @(9,$noTestsCollected))) {
    $r = $pair[0]
    $val = $pair[1]
    $cell = $ws.Cells.Item($r,2)
    $cell.Font.Name = "Arial Unicode MS"
    $cell.Font.Size = 12
    $cell.Font.Color = 0
    $cell.VerticalAlignment = -4108
    $cell.Value2 = $val
}

$errCell = $ws.Cells.Item(7,2)
$errCell.Font.Name = "Consolas"
$errCell.Font.Size = 10
$errCell.Font.Color = 15461355
$errCell.Value2 = $errorOccurred

# ---------------------------------------------------------------------------
# 5. Sheet view: scroll so row 9 is the top row and select B9.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("B9").Select()
$excel.ActiveWindow.ScrollRow = 9
